$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.493404
$ws.Range("H2").Value = 46.48021199999999
$ws.Range("I2").Value = 0.05356331879335558
$ws.Range("J2").Value = 0.05356331879335557
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 21.239540866372
$ws.Range("R2").Value = 191.155867797348
$ws.Range("S2").Value = 0.0005908373168739309
$ws.Range("T2").Value = 0.0005908373168739307

$ws.Range("G3").Value = 15.493404
$ws.Range("H3").Value = 46.48021199999999
$ws.Range("I3").Value = 0.05356331879335558
$ws.Range("J3").Value = 0.05356331879335557
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 1439.549868240964
$ws.Range("R3").Value = 12955.94881416867
$ws.Range("S3").Value = 0.04004511147434213
$ws.Range("T3").Value = 0.04004511147434212

$ws.Range("G4").Value = 15.493404
$ws.Range("H4").Value = 46.48021199999999
$ws.Range("I4").Value = 0.05356331879335558
$ws.Range("J4").Value = 0.05356331879335557
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 460.4404788795359
$ws.Range("R4").Value = 4143.964309915823
$ws.Range("S4").Value = 0.01280844152107144
$ws.Range("T4").Value = 0.01280844152107144

$ws.Range("G5").Value = 15.493404
$ws.Range("H5").Value = 46.48021199999999
$ws.Range("I5").Value = 0.05356331879335558
$ws.Range("J5").Value = 0.05356331879335557
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 4.275265393163999
$ws.Range("R5").Value = 38.47738853847599
$ws.Range("S5").Value = 0.0001189284810680779
$ws.Range("T5").Value = 0.0001189284810680778

$ws.Range("I6").Value = 0.4524333485785276
$ws.Range("J6").Value = 0.4524333485785275
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 179.404055105622
$ws.Range("R6").Value = 1614.636495950598
$ws.Range("S6").Value = 0.004990626267384779
$ws.Range("T6").Value = 0.004990626267384777

$ws.Range("I7").Value = 0.4524333485785276
$ws.Range("J7").Value = 0.4524333485785275
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.338249090733797
$ws.Range("T7").Value = 0.3382490907337969

$ws.Range("I8").Value = 0.4524333485785276
$ws.Range("J8").Value = 0.4524333485785275
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 3889.203140758536
$ws.Range("R8").Value = 35002.82826682682
$ws.Range("S8").Value = 0.1081890782348135
$ws.Range("T8").Value = 0.1081890782348135

$ws.Range("I9").Value = 0.4524333485785276
$ws.Range("J9").Value = 0.4524333485785275
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 36.111889282914
$ws.Range("R9").Value = 325.007003546226
$ws.Range("S9").Value = 0.0010045533425323
$ws.Range("T9").Value = 0.001004553342532299

$ws.Range("G10").Value = 66.835223
$ws.Range("H10").Value = 200.505669
$ws.Range("I10").Value = 0.2310606730563543
$ws.Range("J10").Value = 0.2310606730563542
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 91.62282544375567
$ws.Range("R10").Value = 824.605428993801
$ws.Range("S10").Value = 0.002548745506797011
$ws.Range("T10").Value = 0.00254874550679701

$ws.Range("G11").Value = 66.835223
$ws.Range("H11").Value = 200.505669
$ws.Range("I11").Value = 0.2310606730563543
$ws.Range("J11").Value = 0.2310606730563542
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 6209.90948557886
$ws.Range("R11").Value = 55889.18537020974
$ws.Range("S11").Value = 0.1727460250470145
$ws.Range("T11").Value = 0.1727460250470145

$ws.Range("G12").Value = 66.835223
$ws.Range("H12").Value = 200.505669
$ws.Range("I12").Value = 0.2310606730563543
$ws.Range("J12").Value = 0.2310606730563542
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 1986.241505361932
$ws.Range("R12").Value = 17876.17354825739
$ws.Range("S12").Value = 0.05525287053402012
$ws.Range("T12").Value = 0.05525287053402009

$ws.Range("G13").Value = 66.835223
$ws.Range("H13").Value = 200.505669
$ws.Range("I13").Value = 0.2310606730563543
$ws.Range("J13").Value = 0.2310606730563542
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 18.442578269843
$ws.Range("R13").Value = 165.983204428587
$ws.Range("S13").Value = 0.0005130319685226219
$ws.Range("T13").Value = 0.0005130319685226218

$ws.Range("G14").Value = 76.057215
$ws.Range("H14").Value = 228.171645
$ws.Range("I14").Value = 0.2629426595717627
$ws.Range("J14").Value = 0.2629426595717627
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 104.265036022745
$ws.Range("R14").Value = 938.3853242047051
$ws.Range("S14").Value = 0.002900424002336975
$ws.Range("T14").Value = 0.002900424002336975

$ws.Range("G15").Value = 76.057215
$ws.Range("H15").Value = 228.171645
$ws.Range("I15").Value = 0.2629426595717627
$ws.Range("J15").Value = 0.2629426595717627
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 7066.759108071065
$ws.Range("R15").Value = 63600.83197263959
$ws.Range("S15").Value = 0.196581697159837
$ws.Range("T15").Value = 0.1965816971598369

$ws.Range("G16").Value = 76.057215
$ws.Range("H16").Value = 228.171645
$ws.Range("I16").Value = 0.2629426595717627
$ws.Range("J16").Value = 0.2629426595717627
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 2260.30512706206
$ws.Range("R16").Value = 20342.74614355854
$ws.Range("S16").Value = 0.06287671776861031
$ws.Range("T16").Value = 0.0628767177686103

$ws.Range("G17").Value = 76.057215
$ws.Range("H17").Value = 228.171645
$ws.Range("I17").Value = 0.2629426595717627
$ws.Range("J17").Value = 0.2629426595717627
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 20.987303964315
$ws.Range("R17").Value = 188.885735678835
$ws.Range("S17").Value = 0.0005838206409784596
$ws.Range("T17").Value = 0.0005838206409784595
